# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "30.367.55"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.871.32"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5: BNB
$cell = $ws.Range("D5")
$savedStyle = $cell.Style
$cell.Value = "'244.68"
$cell.Style = $savedStyle
$ws.Range("E5").Value = "  -1.95%  "

# Row 6: USDC
$ws.Range("E6").Value = "  -0.03%  "

# Row 7: XRP
$cell = $ws.Range("D7")
$savedStyle = $cell.Style
$cell.Value = "'0.4724"
$cell.Style = $savedStyle
$ws.Range("E7").Value = "  -0.56%  "

# Row 8: Cardano
$cell = $ws.Range("D8")
$savedStyle = $cell.Style
$cell.Value = "'0.2867"
$cell.Style = $savedStyle
$ws.Range("E8").Value = "  -2.25%  "

# Row 9: Dogecoin
$cell = $ws.Range("D9")
$savedStyle = $cell.Style
$cell.Value = "'0.06484"
$cell.Style = $savedStyle
$ws.Range("E9").Value = "  -0.82%  "

# Row 10: Solana
$cell = $ws.Range("D10")
$savedStyle = $cell.Style
$cell.Value = "'21.76"
$cell.Style = $savedStyle
$ws.Range("E10").Value = "  -1.05%  "

# Row 11: Litecoin
$cell = $ws.Range("D11")
$savedStyle = $cell.Style
$cell.Value = "'100.13"
$cell.Style = $savedStyle
$ws.Range("E11").Value = "  +2.99%  "

# Row 12: TRON
$cell = $ws.Range("D12")
$savedStyle = $cell.Style
$cell.Value = "'0.07796"
$cell.Style = $savedStyle
$ws.Range("E12").Value = "  -0.01%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.871.57"
$ws.Range("E13").Value = "  -1.04%  "

# Row 14: Polygon
$cell = $ws.Range("D14")
$savedStyle = $cell.Style
$cell.Value = "'0.7272"
$cell.Style = $savedStyle
$ws.Range("E14").Value = "  -1.56%  "

# Row 15: Polkadot
$cell = $ws.Range("D15")
$savedStyle = $cell.Style
$cell.Value = "'5.163"
$cell.Style = $savedStyle
$ws.Range("E15").Value = "  -1.53%  "

# Row 16: BitcoinCash
$cell = $ws.Range("D16")
$savedStyle = $cell.Style
$cell.Value = "'283.49"
$cell.Style = $savedStyle
$ws.Range("E16").Value = "  -0.45%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "30.352.68"
$ws.Range("E17").Value = "  -1.52%  "

# Row 18: Avalanche
$cell = $ws.Range("D18")
$savedStyle = $cell.Style
$cell.Value = "'13.08"
$cell.Style = $savedStyle
$ws.Range("E18").Value = "  -0.96%  "

# Row 19: Dai
$cell = $ws.Range("D19")
$savedStyle = $cell.Style
$cell.Value = "'1.000"
$cell.Style = $savedStyle
$ws.Range("E19").Value = "  -0.01%  "

# Row 20: ShibaInu
$cell = $ws.Range("D20")
$savedStyle = $cell.Style
$cell.Value = "'0.000007469"
$cell.Style = $savedStyle
$ws.Range("E20").Value = "  -1.37%  "

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.115.58"
$ws.Range("E21").Value = "  -1.19%  "

# Row 22: Uniswap
$cell = $ws.Range("D22")
$savedStyle = $cell.Style
$cell.Value = "'5.324"
$cell.Style = $savedStyle
$ws.Range("E22").Value = "  -0.18%  "

# Row 23: BinanceUSD
$cell = $ws.Range("D23")
$savedStyle = $cell.Style
$cell.Value = "'0.9999"
$cell.Style = $savedStyle
$ws.Range("E23").Value = "  -0.10%  "

# Row 24: Chainlink
$cell = $ws.Range("D24")
$savedStyle = $cell.Style
$cell.Value = "'6.323"
$cell.Style = $savedStyle
$ws.Range("E24").Value = "  +0.90%  "

# Row 25: Monero
$cell = $ws.Range("D25")
$savedStyle = $cell.Style
$cell.Value = "'162.87"
$cell.Style = $savedStyle
$ws.Range("E25").Value = "  -0.97%  "

# Row 26: Cosmos
$cell = $ws.Range("D26")
$savedStyle = $cell.Style
$cell.Value = "'9.026"
$cell.Style = $savedStyle
$ws.Range("E26").Value = "  -2.34%  "

# Row 27: EthereumClassic
$cell = $ws.Range("D27")
$savedStyle = $cell.Style
$cell.Value = "'18.92"
$cell.Style = $savedStyle
$ws.Range("E27").Value = "  -0.29%  "

# Row 28: LidoDAOToken
$cell = $ws.Range("D28")
$savedStyle = $cell.Style
$cell.Value = "'1.890"
$cell.Style = $savedStyle
$ws.Range("E28").Value = "  -1.79%  "

# Row 29: Stellar
$cell = $ws.Range("D29")
$savedStyle = $cell.Style
$cell.Value = "'0.09651"
$cell.Style = $savedStyle
$ws.Range("E29").Value = "  -0.84%  "

# Row 30: Toncoin
$cell = $ws.Range("D30")
$savedStyle = $cell.Style
$cell.Value = "'1.322"
$cell.Style = $savedStyle
$ws.Range("E30").Value = "  -1.53%  "

# Row 31: PancakeSwap
$cell = $ws.Range("D31")
$savedStyle = $cell.Style
$cell.Value = "'1.489"
$cell.Style = $savedStyle
$ws.Range("E31").Value = "  -0.79%  "

# Row 32: Filecoin
$cell = $ws.Range("D32")
$savedStyle = $cell.Style
$cell.Value = "'4.222"
$cell.Style = $savedStyle
$ws.Range("E32").Value = "  -1.98%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.98%  "

# Row 34: Hedera
$cell = $ws.Range("D34")
$savedStyle = $cell.Style
$cell.Value = "'0.04802"
$cell.Style = $savedStyle
$ws.Range("E34").Value = "  -1.91%  "

# Row 35: ARBITRUM
$cell = $ws.Range("D35")
$savedStyle = $cell.Style
$cell.Value = "'1.124"
$cell.Style = $savedStyle
$ws.Range("E35").Value = "  -0.51%  "

# Row 36: ImmutableX
$cell = $ws.Range("D36")
$savedStyle = $cell.Style
$cell.Value = "'0.6885"
$cell.Style = $savedStyle
$ws.Range("E36").Value = "  -1.66%  "

# Row 37: HuobiToken
$cell = $ws.Range("D37")
$savedStyle = $cell.Style
$cell.Value = "'2.720"
$cell.Style = $savedStyle
$ws.Range("E37").Value = "  -0.14%  "

# Row 38: VeChain
$cell = $ws.Range("D38")
$savedStyle = $cell.Style
$cell.Value = "'0.01898"
$cell.Style = $savedStyle
$ws.Range("E38").Value = "  -0.74%  "

# Row 39: MXToken
$cell = $ws.Range("D39")
$savedStyle = $cell.Style
$cell.Value = "'2.838"
$cell.Style = $savedStyle
$ws.Range("E39").Value = "  +0.80%  "

# Row 40: Aave
$cell = $ws.Range("D40")
$savedStyle = $cell.Style
$cell.Value = "'76.38"
$cell.Style = $savedStyle
$ws.Range("E40").Value = "  -0.26%  "

# Row 41: FraxShare
$cell = $ws.Range("D41")
$savedStyle = $cell.Style
$cell.Value = "'6.282"
$cell.Style = $savedStyle
$ws.Range("E41").Value = "  -1.32%  "

# Row 42: RenderToken
$cell = $ws.Range("D42")
$savedStyle = $cell.Style
$cell.Value = "'1.950"
$cell.Style = $savedStyle
$ws.Range("E42").Value = "  -3.49%  "

# Row 43: TheSandbox
$cell = $ws.Range("D43")
$savedStyle = $cell.Style
$cell.Value = "'0.4210"
$cell.Style = $savedStyle
$ws.Range("E43").Value = "  -1.39%  "

# Row 44: PaxDollar
$cell = $ws.Range("D44")
$savedStyle = $cell.Style
$cell.Value = "'1.000"
$cell.Style = $savedStyle
$ws.Range("E44").Value = "  -0.04%  "

# Row 45: TrustWalletToken
$cell = $ws.Range("D45")
$savedStyle = $cell.Style
$cell.Value = "'0.8220"
$cell.Style = $savedStyle
$ws.Range("E45").Value = "  -1.77%  "

# Row 46: Quant
$cell = $ws.Range("D46")
$savedStyle = $cell.Style
$cell.Value = "'100.87"
$cell.Style = $savedStyle
$ws.Range("E46").Value = "  -0.73%  "

# Row 47: EnergySwap
$cell = $ws.Range("D47")
$savedStyle = $cell.Style
$cell.Value = "'9.729"
$cell.Style = $savedStyle
$ws.Range("E47").Value = "  +2.61%  "

# Row 48: Aptos
$cell = $ws.Range("D48")
$savedStyle = $cell.Style
$cell.Value = "'7.007"
$cell.Style = $savedStyle
$ws.Range("E48").Value = "  -1.00%  "

# Row 49: Elrond
$cell = $ws.Range("D49")
$savedStyle = $cell.Style
$cell.Value = "'34.87"
$cell.Style = $savedStyle
$ws.Range("E49").Value = "  -2.45%  "

# Row 50: Cronos
$cell = $ws.Range("D50")
$savedStyle = $cell.Style
$cell.Value = "'0.05756"
$cell.Style = $savedStyle
$ws.Range("E50").Value = "  -0.14%  "

# Row 51: Maker
$cell = $ws.Range("D51")
$savedStyle = $cell.Style
$cell.Value = "'882.26"
$cell.Style = $savedStyle
$ws.Range("E51").Value = "  -4.52%  "
